$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.6773685352888915
$ws.Range("E2").Value = 0.6773685352888915

# Row 3
$ws.Range("D3").Value = [double]"2.900456263355253E-05"
$ws.Range("E3").Value = [double]"2.900456263355253E-05"

# Row 4
$ws.Range("D4").Value = 0.002590123258075552
$ws.Range("E4").Value = 0.002590123258075552

# Row 5
$ws.Range("D5").Value = 0.0007905143580800866
$ws.Range("E5").Value = 0.0007905143580800866

# Row 6
$ws.Range("D6").Value = 0.8074772992877954
$ws.Range("E6").Value = 0.8074772992877954

# Row 7
$ws.Range("D7").Value = 0.7288369672418282
$ws.Range("E7").Value = 0.2711630327581718

# Row 8
$ws.Range("C8").Value = $false
$ws.Range("D8").Value = [double]"1.202827873786564E-06"
$ws.Range("E8").Value = 0.9999987971721263

# Row 9
$ws.Range("C9").Value = $false
$ws.Range("D9").Value = [double]"3.602875892714194E-12"
$ws.Range("E9").Value = 0.9999999999963971

# Row 10
$ws.Range("C10").Value = $false
$ws.Range("D10").Value = [double]"2.771316706537062E-06"
$ws.Range("E10").Value = 0.9999972286832934

# Row 11
$ws.Range("C11").Value = $false
$ws.Range("D11").Value = 0.1771414761422895
$ws.Range("E11").Value = 0.8228585238577105
$ws.Range("F11").Value = 5.760561943054199
$ws.Range("G11").Value = 0.4

# Row 12
$ws.Range("D12").Value = 0.8537451610423439
$ws.Range("E12").Value = 0.8537451610423439

# Row 13
$ws.Range("D13").Value = [double]"8.361315539162392E-06"
$ws.Range("E13").Value = [double]"8.361315539162392E-06"

# Row 14
$ws.Range("D14").Value = 0.0003712103656480031
$ws.Range("E14").Value = 0.0003712103656480031

# Row 15
$ws.Range("D15").Value = [double]"8.346299101834615E-06"
$ws.Range("E15").Value = [double]"8.346299101834615E-06"

# Row 16
$ws.Range("D16").Value = 0.910891109550215
$ws.Range("E16").Value = 0.910891109550215

# Row 17
$ws.Range("D17").Value = 0.8591180471349766
$ws.Range("E17").Value = 0.1408819528650234

# Row 18
$ws.Range("C18").Value = $false
$ws.Range("D18").Value = [double]"2.182498765102871E-06"
$ws.Range("E18").Value = 0.9999978175012348

# Row 19
$ws.Range("C19").Value = $false
$ws.Range("D19").Value = [double]"5.828603003825485E-19"
$ws.Range("E19").Value = 1

# Row 20
$ws.Range("C20").Value = $false
$ws.Range("D20").Value = [double]"1.989691771583631E-07"
$ws.Range("E20").Value = 0.9999998010308229

# Row 21
$ws.Range("C21").Value = $false
$ws.Range("D21").Value = 0.04089707394420425
$ws.Range("E21").Value = 0.9591029260557957
$ws.Range("F21").Value = 7.814073085784912
$ws.Range("G21").Value = 0.4

$wb.Save()
